$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.909.66"
$ws.Range("E2").Value = "  -0.41%  "
$ws.Range("D3").Value = "1.668.35"
$ws.Range("E3").Value = "  +1.05%  "
$ws.Range("E4").Value = "  +0.11%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "215.48"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.21%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.523"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +2.92%  "
$ws.Range("E7").Value = "  +0.07%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.0620"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +1.15%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.249"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -0.01%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "20.26"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +2.70%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0894"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +3.31%  "
$ws.Range("D12").Value = "1.904.16"
$ws.Range("E12").Value = "  +1.09%  "
$ws.Range("D13").Value = "1.698.71"
$ws.Range("E13").Value = "  +2.86%  "
$ws.Range("E14").Value = "  +0.32%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.526"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +1.39%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "66.08"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +1.46%  "
$ws.Range("D17").Value = "26.924.78"
$ws.Range("E17").Value = "  -0.29%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "233.84"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -2.01%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "7.99"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +1.80%  "
$ws.Range("E20").Value = "  +0.40%  "
$ws.Range("E21").Value = "  +0.07%  "
$ws.Range("E22").Value = "  -0.14%  "
$ws.Range("E23").Value = "  -1.12%  "
$ws.Range("E24").Value = "  -3.16%  "
$ws.Range("E25").Value = "  +0.46%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "7.14"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.59%  "
$ws.Range("E27").Value = "  +0.32%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "15.89"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +0.56%  "
$ws.Range("E29").Value = "  +0.11%  "
$ws.Range("E30").Value = "  +0.07%  "
$ws.Range("E31").Value = "  +0.28%  "
$ws.Range("E32").Value = "  +1.71%  "
$ws.Range("D33").Value = "1.455.18"
$ws.Range("E33").Value = "  -3.78%  "
$ws.Range("E34").Value = "  +2.08%  "
$ws.Range("E35").Value = "  +3.32%  "
$ws.Range("E36").Value = "  -0.15%  "
$ws.Range("E37").Value = "  +0.81%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.901"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +1.77%  "
$ws.Range("E39").Value = "  +0.86%  "
$ws.Range("E40").Value = "  -3.47%  "
$ws.Range("E41").Value = "  +0.08%  "
$ws.Range("E42").Value = "  +1.07%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "66.14"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +0.24%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.977"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +6.77%  "
$ws.Range("D45").Value = "1.815.72"
$ws.Range("E45").Value = "  +1.40%  "
$ws.Range("E46").Value = "  +1.37%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "90.70"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +1.29%  "
$ws.Range("E48").Value = "  +1.28%  "
$ws.Range("E49").Value = "  -0.93%  "
$ws.Range("E50").Value = "  +4.29%  "
$ws.Range("E51").Value = "  +0.00%  "
